# Daily attendance processing - reorder "Recorded By" (column G) entries.
#
# Rule observed in the target diff: within each comma-separated "Recorded By"
# list, if the LAST entry is an email address (contains "@"), it is rotated
# to the front of the list (moved from the end to the start), with the
# remaining entries keeping their relative order. Lists whose last entry is
# not an email (e.g. "System") - or that only have a single entry - are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G = "Recorded By"
$col = 7

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    if (-not $text.Contains(",")) {
        continue
    }

    $parts = @($text.Split(",") | ForEach-Object { $_.Trim() })

    if ($parts.Count -le 1) {
        continue
    }

    $lastIdx = $parts.Count - 1
    $lastPart = $parts[$lastIdx]

    if ($lastPart.Contains("@")) {
        $newParts = @($lastPart) + $parts[0..($lastIdx - 1)]
        $newText = $newParts -join ", "
        if ($newText -ne $text) {
            $cell.Value = $newText
        }
    }
}
